$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row above row 30. This shifts rows 30:57 down to 31:58,
# carrying along all their existing values/formatting, matching the
# append of a new (most-recent) September notification entry at the
# top of the reverse-chronological R/S (September_Details/September_Date) list.
$ws.Rows("30:30").Insert()

# Populate the newly inserted row 30 with the new notification entry.
$ws.Range("R30").Value = "broker"
$ws.Range("S30").Value = "2024-09-04 21:20:47"
